# QA Round 2: deep quality optimization - compliance, diversification, UX improvements
# Co-authored-by: Cursor <cursoragent@cursor.com>
#
# This script:
#  1) Updates a handful of script lines on the "MaddisonJourney" sheet.
#  2) Splits the old "cumcontrol" sheet into "cumcontrol1" (edited in place)
#     and "cumcontrol2" (built from the old "dickpic" sheet content, edited),
#     re-creating a fresh, unmodified "dickpic" sheet right after it so the
#     tab order stays: ... done2, cumcontrol1, cumcontrol2, dickpic, boosters.
#  3) Updates the EDGE/DELAY/SYNC/CONTROL copy across both new sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) MaddisonJourney script tweaks
# ---------------------------------------------------------------------------
$journey = $wb.Worksheets.Item("MaddisonJourney")
$journey.Range("B4").Value = "finish with me Daddy"
$journey.Range("B5").Value = "don't go anywhere"
$journey.Range("B11").Value = "I can't control myself anymore Daddy"
$journey.Range("B14").Value = "fuuuck"
$journey.Range("B22").Value = "gimme a minute"

# ---------------------------------------------------------------------------
# 2) Rename "cumcontrol" -> "cumcontrol1"
# ---------------------------------------------------------------------------
$cumcontrol1 = $wb.Worksheets.Item("cumcontrol")
$cumcontrol1.Name = "cumcontrol1"

# ---------------------------------------------------------------------------
# 3) Duplicate "dickpic" so we keep an untouched copy, then repurpose the
#    original "dickpic" sheet (renamed to "cumcontrol2") for the new,
#    edited control-themed variants. Placing the duplicate right after the
#    original preserves tab order: cumcontrol1, cumcontrol2(old dickpic), dickpic(new copy), boosters
# ---------------------------------------------------------------------------
$dickpicOriginal = $wb.Worksheets.Item("dickpic")
$dickpicOriginal.Copy($null, $dickpicOriginal)
$dickpicCopy = $wb.ActiveSheet

$cumcontrol2 = $dickpicOriginal
$cumcontrol2.Name = "cumcontrol2"
$dickpicCopy.Name = "dickpic"

# ---------------------------------------------------------------------------
# 4) cumcontrol1 text edits (rows 2-7)
# ---------------------------------------------------------------------------
$cumcontrol1.Range("B2").Value = "trust me you want to edge just a little longer for this one"

$cumcontrol1.Range("B3").Value = "you're not done until I say you are... open this"
$cumcontrol1.Range("C3").Value = "DELAY. Send PPV."

$cumcontrol1.Range("B4").Value = "I'm right there too Daddy, let's finish this... but you need to see this first"
$cumcontrol1.Range("C4").Value = "SYNC variant. Send PPV."

$cumcontrol1.Range("B5").Value = "now... right now, with me. open this"
$cumcontrol1.Range("C5").Value = "SYNC. Send PPV."

$cumcontrol1.Range("B6").Value = "not a chance... you're going to wait until I say so Daddy"

$cumcontrol1.Range("B7").Value = "I didn't say you could cum yet Daddy"
$cumcontrol1.Range("C7").Value = "CONTROL."

# ---------------------------------------------------------------------------
# 5) cumcontrol2 text edits (rows 2-7) - names (A), text (B), notes (C)
# ---------------------------------------------------------------------------
$cumcontrol2.Range("A2").Value = "delay2"
$cumcontrol2.Range("B2").Value = "edge for me Daddy... just a little more... this last one is everything"
$cumcontrol2.Range("C2").Value = "DELAY variant."

$cumcontrol2.Range("A3").Value = "delay1"
$cumcontrol2.Range("B3").Value = "hold it... what I'm about to send is the best one and you'll want to last for it"
$cumcontrol2.Range("C3").Value = "DELAY. Send PPV."

$cumcontrol2.Range("A4").Value = "sync2"
$cumcontrol2.Range("B4").Value = "I want us to finish at the same time Daddy... this one will push you over"
$cumcontrol2.Range("C4").Value = "SYNC variant."

$cumcontrol2.Range("A5").Value = "sync1"
$cumcontrol2.Range("B5").Value = "okay you earned it... let's go together, open this"
$cumcontrol2.Range("C5").Value = "SYNC. Send PPV."

$cumcontrol2.Range("A6").Value = "edge2"
$cumcontrol2.Range("B6").Value = "if you finish without my permission I'll be annoyed"
$cumcontrol2.Range("C6").Value = "EDGE variant."

$cumcontrol2.Range("A7").Value = "edge1"
$cumcontrol2.Range("B7").Value = "slow down Daddy, I'm in control here"
$cumcontrol2.Range("C7").Value = "CONTROL."
